# Updates cryptos list price/volume data (and a few reordered rows) to match
# the latest scrape, as produced by the scheduled GitHub Actions refresh job.
# Note: numeric-looking Price strings (column D) are prefixed with a leading
# apostrophe so Excel stores them as literal text (preserving formats like
# "1.00" or "0.522") instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.654.80'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '1.626.97'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').Value = '''212.22'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '''0.522'
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').Value = '''22.89'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = '''0.261'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '1.861.19'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '1.640.66'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('D16').Value = '''64.45'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '27.647.92'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '''229.50'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '0.0₃0722'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '''7.57'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '''9.95'
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('E24').Value = '  +6.47%  '
$ws.Range('D25').Value = '''149.75'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '''6.89'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').Value = '''1.01'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '''0.111'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('D29').Value = '''15.57'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('D31').Value = '''0.0483'
$ws.Range('E31').Value = '  -1.14%  '
$ws.Range('D32').Value = '''3.29'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').Value = '1.462.34'
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('D37').Value = '''0.567'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').Value = '''0.869'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').Value = '''0.913'
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('D41').Value = '''69.32'
$ws.Range('E41').Value = '  +6.71%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.01'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '''1.02'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''5.43'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '''2.23'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').Value = '1.771.52'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').Value = '''1.70'
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('D49').Value = '''85.84'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0999'
$ws.Range('E50').Value = '  -6.75%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0986'
$ws.Range('E51').Value = '  -0.21%  '
